$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169334530830383
$ws.Range("B1").Value = 2.44109320640564
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.362691402435303
$ws.Range("E1").Value = 1.237056136131287
